# Extend the capabilities of Data driven approach
# - Rename TC1 -> LoginPage_TC1 on both sheets (shared string updated in place)
# - Add login credentials to the "Test Data" sheet (B2/C2), entered with a
#   leading apostrophe so Excel stores them as text (quotePrefix) and
#   center-aligned like the existing data cell.
# - Add the "loginToApp" keyword to the "Business Flow" sheet (B2)
# - Make "Business Flow" the active/selected sheet (was "Test Data")
# - Update the remembered cell selection on both sheets
# - Widen column A on both sheets to fit the longer TC id text

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Business Flow"
$ws2 = $wb.Worksheets.Item(2)   # "Test Data"

# --- Cell content -----------------------------------------------------

# Rename the shared "TC1" value to "LoginPage_TC1" on both sheets so the
# shared string itself is updated instead of creating a brand-new one.
$ws1.Range("A2").Value() = "LoginPage_TC1"
$ws2.Range("A2").Value() = "LoginPage_TC1"

# New login credentials on the "Test Data" sheet, centered like A2, stored
# as text (leading apostrophe forces the quote-prefix / text style).
$ws2.Range("B2").Value() = "'tejeshkumar.gangari@gmail.com"
$ws2.Range("C2").Value() = "'Test@123"
$ws2.Range("B2:C2").HorizontalAlignment = -4108

# New keyword on the "Business Flow" sheet.
$ws1.Range("B2").Value() = "loginToApp"

# --- Column widths ------------------------------------------------------

$ws1.Columns.Item(1).ColumnWidth = 13.501302083333334
$ws2.Columns.Item(1).ColumnWidth = 13.501302083333334
$ws2.Columns.Item(2).ColumnWidth = 31.16796875
$ws2.Columns.Item(3).ColumnWidth = 11.66796875

# --- Sheet view / selection ---------------------------------------------

# "Test Data" used to be the active tab; flip it to "Business Flow" and
# move each sheet's remembered selection.
$ws2.Range("F8").Select()
$ws1.Activate()
$ws1.Range("D7").Select()
